# Applies the cell-value corrections captured in the upstream diff for
# "Kujata_Profits.xlsx" (the Leve-profit market-board refresh). Every
# touched cell is a literal cached value (no formulas in this workbook),
# so each change is expressed as a direct Range.Value assignment; cells
# that the diff removes entirely are cleared so they go back to "no cell".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 442.65518
$ws.Range("I33").Value = 420.22726
$ws.Range("K33").Value = 420.22726
$ws.Range("M33").Value = -191.22726
# Row 62
$ws.Range("H62").Value = 18521736
$ws.Range("I62").Value = 22225482
$ws.Range("J62").Value = 3006
$ws.Range("K62").Value = 22225482
$ws.Range("L62").Value = 3006
$ws.Range("M62").Value = -22224858
$ws.Range("N62").Value = -4254
# Row 65
$ws.Range("H65").Value = 18521736
$ws.Range("I65").Value = 22225482
$ws.Range("J65").Value = 3006
$ws.Range("K65").Value = 111127410
$ws.Range("L65").Value = 15030
$ws.Range("M65").Value = -111124290
$ws.Range("N65").Value = -21270
# Row 138
$ws.Range("H138").Value = 1620.11
$ws.Range("I138").Value = 927.6087
$ws.Range("J138").Value = 1826.961
$ws.Range("K138").Value = 2782.8261
$ws.Range("L138").Value = 5480.883
$ws.Range("M138").Value = 2357.1739
$ws.Range("N138").Value = -15760.883

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4784.327
$ws.Range("I32").Value = 4575.76
$ws.Range("K32").Value = 4575.76
$ws.Range("M32").Value = -4288.76
# Row 61
$ws.Range("H61").Value = 37037930
$ws.Range("I61").Value = 40000812
$ws.Range("J61").Value = 1890
$ws.Range("K61").Value = 40000812
$ws.Range("L61").Value = 1890
$ws.Range("M61").Value = -40000600
$ws.Range("N61").Value = -2314
# Row 136
$ws.Range("H136").Value = 37037930
$ws.Range("I136").Value = 40000812
$ws.Range("J136").Value = 1890
$ws.Range("K136").Value = 120002436
$ws.Range("L136").Value = 5670
$ws.Range("M136").Value = -119999886
$ws.Range("N136").Value = -10770

$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 94
$ws.Range("H94").Value = 22727978
$ws.Range("I94").Value = 25000656
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 25000656
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -25000205
$ws.Range("N94").Value = -2102

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 772.7273
$ws.Range("I105").Value = 766.6667
$ws.Range("K105").Value = 766.6667
$ws.Range("M105").Value = 980.3333
# Row 132
$ws.Range("H132").Value = 3412.6667
$ws.Range("I132").Value = 3161.5
$ws.Range("K132").Value = 9484.5
$ws.Range("M132").Value = -6954.5
# Row 137
$ws.Range("H137").Value = 62890
$ws.Range("J137").Value = 65468
$ws.Range("L137").Value = 65468
$ws.Range("N137").Value = -75668

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 435.875
$ws.Range("I5").Value = 406.7143
$ws.Range("K5").Value = 1220.1429
$ws.Range("M5").Value = -1108.1429
# Row 114
$ws.Range("H114").Value = 509.66666
$ws.Range("I114").Value = 390
$ws.Range("J114").Value = 646.4286
$ws.Range("K114").Value = 1170
$ws.Range("L114").Value = 1939.2858
$ws.Range("M114").Value = 2084
$ws.Range("N114").Value = -8447.2858
# Row 135
$ws.Range("H135").Value = 435.875
$ws.Range("I135").Value = 406.7143
$ws.Range("K135").Value = 3660.4287
$ws.Range("M135").Value = -1125.4287

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 2333.3333
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -712
$ws.Range("N19").Value = -5576
# Row 20
$ws.Range("H20").Value = 33335466
$ws.Range("I20").Value = 50000000
$ws.Range("J20").Value = 6400
$ws.Range("K20").Value = 50000000
$ws.Range("L20").Value = 6400
$ws.Range("M20").Value = -49999755
$ws.Range("N20").Value = -6890

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1536.4286
$ws.Range("I7").Value = 1633.3334
$ws.Range("J7").Value = 1463.75
$ws.Range("K7").Value = 1633.3334
$ws.Range("L7").Value = 1463.75
$ws.Range("M7").Value = -1521.3334
$ws.Range("N7").Value = -1687.75
# Row 40
$ws.Range("H40").Value = 3095.1052
$ws.Range("J40").Value = 5617.3335
$ws.Range("L40").Value = 5617.3335
$ws.Range("N40").Value = -5889.3335
# Row 46
$ws.Range("H46").Value = 6424.143
$ws.Range("I46").Value = 984.5
$ws.Range("K46").Value = 984.5
$ws.Range("M46").Value = -796.5
# Row 68
$ws.Range("H68").Value = 1823
$ws.Range("I68").Value = 1808.25
$ws.Range("K68").Value = 1808.25
$ws.Range("M68").Value = -1059.25
# Row 71
$ws.Range("H71").Value = 1823
$ws.Range("I71").Value = 1808.25
$ws.Range("K71").Value = 9041.25
$ws.Range("M71").Value = -5297.25
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 126
$ws.Range("H126").Value = 1536.4286
$ws.Range("I126").Value = 1633.3334
$ws.Range("J126").Value = 1463.75
$ws.Range("K126").Value = 4900.0002
$ws.Range("L126").Value = 4391.25
$ws.Range("M126").Value = -2430.0002
$ws.Range("N126").Value = -9331.25
# Row 133
$ws.Range("H133").Value = 46663
$ws.Range("J133").Value = 46663
$ws.Range("L133").Value = 46663
$ws.Range("N133").Value = -51723

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 80004
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 80004
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 80004
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -80232
# Row 5
$ws.Range("H5").Value = 30333334
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 30333334
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 30333334
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -30333558
# Row 7
$ws.Range("H7").Value = 2931191.8
$ws.Range("I7").Value = 195
$ws.Range("K7").Value = 195
$ws.Range("M7").Value = -82
# Row 18
$ws.Range("H18").Value = 407
$ws.Range("J18").Value = 407
$ws.Range("L18").Value = 407
$ws.Range("N18").Value = -753
# Row 38
$ws.Range("H38").Value = 1525
$ws.Range("I38").Value = 1525
$ws.Range("K38").Value = 1525
$ws.Range("M38").Value = -1052
# Row 62
$ws.Range("H62").Value = 166676670
$ws.Range("I62").Value = 250005000
$ws.Range("J62").Value = 20003
$ws.Range("K62").Value = 250005000
$ws.Range("L62").Value = 20003
$ws.Range("M62").Value = -250004376
$ws.Range("N62").Value = -21251
# Row 65
$ws.Range("H65").Value = 166676670
$ws.Range("I65").Value = 250005000
$ws.Range("J65").Value = 20003
$ws.Range("K65").Value = 1250025000
$ws.Range("L65").Value = 100015
$ws.Range("M65").Value = -1250021880
$ws.Range("N65").Value = -106255
# Row 104
$ws.Range("H104").Value = 19450
$ws.Range("J104").Value = 19450
$ws.Range("L104").Value = 19450
$ws.Range("N104").Value = -26438
# Row 126
$ws.Range("H126").Value = 45455360
$ws.Range("I126").Value = 52632210
$ws.Range("K126").Value = 157896630
$ws.Range("M126").Value = -157894160

